$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is forced to Text format before assignment so that
# numeric-looking strings (e.g. "235.88") are stored as text, matching
# the original inlineStr cell type. The style is then reset back to
# "Normal" so no stray cell style/number-format is left behind.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.652.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.415.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.096.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.669.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.35%  "
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.37%  "
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.481.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -9.64%  "
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "MXToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "FraxShare"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.305.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.29%  "
$ws.Range("E51").Style = "Normal"

Write-Output "Applied all cell updates"
